$d = $word.ActiveDocument

# --------------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits right
#    after the H1 title at the top of the document.
# --------------------------------------------------------------------------
$p2 = $d.Paragraphs(2)
if ($p2.Range.Text -like "Meta description*") {
    $p2.Range.Delete()
}

# --------------------------------------------------------------------------
# 2) Insert a new bold paragraph ("Play Congo King Quad Shot Free Online
#    Review | 100 Paylines") right before the final paragraph (the one that
#    currently holds the DALLE image prompt).
# --------------------------------------------------------------------------
$last = $d.Paragraphs($d.Paragraphs.Count)
$insertPoint = $last.Range.Start
$insertRange = $d.Range($insertPoint, $insertPoint)

$snippet = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' +
           '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Congo King Quad Shot Free Online Review | 100 Paylines</w:t></w:r></w:p>' +
           '</w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

$insertRange.InsertXML($snippet)

# The inserted fragment merges into the paragraph at the insertion point, so
# split it back into two separate paragraphs by adding a paragraph break
# right after the newly-inserted heading text.
$searchRange = $d.Range($insertPoint, $d.Content.End)
$searchRange.Find.Execute("Play Congo King Quad Shot Free Online Review | 100 Paylines", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$breakPoint = $d.Range($searchRange.End, $searchRange.End)
$breakPoint.InsertParagraphAfter()

# --------------------------------------------------------------------------
# 3) Replace the DALLE image-prompt text (now in the true last paragraph)
#    with the meta-description text, keeping its italic formatting intact.
# --------------------------------------------------------------------------
$oldText = 'DALLE, please create a feature image fitting the game "Congo King Quad Shot" that meets the following requirements: - The image should be in cartoon style - The image should feature a happy Maya warrior with glasses. The image should capture the adventurous spirit of the game and convey the excitement of exploring through the jungle. Please ensure the colors used in the image are vibrant and eye-catching. The image should be appealing and encourage players to take a chance on the game.'
$newText = 'Experience the thrill of winning with Congo King Quad Shot. 100 paylines, x4 jackpots, and 15 free spins. Play free online here.'

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
